$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.960525512695312
$ws.Range("C2").Value = 5.683907985687256
$ws.Range("D2").Value = 13.921051979064941
$ws.Range("E2").Value = 57.85714340209961
